# chore: update Sheets via scheduled runner
# Refresh leve-profit derived values (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ,
# LeveProfit NQ/HQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect
# newly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1467.909
$ws.Range("I6").Value = 207.83333
$ws.Range("J6").Value = 2980
$ws.Range("K6").Value = 623.49999
$ws.Range("L6").Value = 8940
$ws.Range("M6").Value = -511.49999
$ws.Range("N6").Value = -9164

$ws.Range("H8").Value = 932.2941
$ws.Range("I8").Value = 104.454544
$ws.Range("J8").Value = 2450
$ws.Range("K8").Value = 313.363632
$ws.Range("L8").Value = 7350
$ws.Range("M8").Value = -174.363632
$ws.Range("N8").Value = -7628

$ws.Range("H18").Value = 2785498
$ws.Range("I18").Value = 4632496
$ws.Range("J18").Value = 15001
$ws.Range("K18").Value = 4632496
$ws.Range("L18").Value = 15001
$ws.Range("M18").Value = -4632212
$ws.Range("N18").Value = -15569

$ws.Range("H116").Value = 2871.718
$ws.Range("I116").Value = 2430.75
$ws.Range("J116").Value = 3335.8948
$ws.Range("K116").Value = 2430.75
$ws.Range("L116").Value = 3335.8948
$ws.Range("M116").Value = 1011.25
$ws.Range("N116").Value = -10219.8948

$ws.Range("H131").Value = 3164.25
$ws.Range("I131").Value = 2627
$ws.Range("J131").Value = 4162
$ws.Range("K131").Value = 7881
$ws.Range("L131").Value = 12486
$ws.Range("M131").Value = -2841
$ws.Range("N131").Value = -22566

$ws.Range("H132").Value = 4083146.2
$ws.Range("I132").Value = 4652153
$ws.Range("J132").Value = 5265
$ws.Range("K132").Value = 13956459
$ws.Range("L132").Value = 15795
$ws.Range("M132").Value = -13953929
$ws.Range("N132").Value = -20855

$ws.Range("H137").Value = 2859946
$ws.Range("I137").Value = 4350840
$ws.Range("J137").Value = 2399.1667
$ws.Range("K137").Value = 13052520
$ws.Range("L137").Value = 7197.500100000001
$ws.Range("M137").Value = -13049970
$ws.Range("N137").Value = -12297.5001

$ws.Range("H138").Value = 3536.6135
$ws.Range("I138").Value = 2144.7666
$ws.Range("J138").Value = 6519.143
$ws.Range("K138").Value = 6434.2998
$ws.Range("L138").Value = 19557.429
$ws.Range("M138").Value = -1294.2998
$ws.Range("N138").Value = -29837.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10001648
$ws.Range("I2").Value = 31251012
$ws.Range("J2").Value = 1947
$ws.Range("K2").Value = 31251012
$ws.Range("L2").Value = 1947
$ws.Range("M2").Value = -31250899
$ws.Range("N2").Value = -2173

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H74").Value = 928.3333
$ws.Range("I74").Value = 837.8333
$ws.Range("J74").Value = 1199.8334
$ws.Range("K74").Value = 837.8333
$ws.Range("L74").Value = 1199.8334
$ws.Range("M74").Value = 36.16669999999999
$ws.Range("N74").Value = -2947.8334

$ws.Range("H77").Value = 928.3333
$ws.Range("I77").Value = 837.8333
$ws.Range("J77").Value = 1199.8334
$ws.Range("K77").Value = 4189.1665
$ws.Range("L77").Value = 5999.166999999999
$ws.Range("M77").Value = 178.8334999999997
$ws.Range("N77").Value = -14735.167

$ws.Range("H116").Value = 10001648
$ws.Range("I116").Value = 31251012
$ws.Range("J116").Value = 1947
$ws.Range("K116").Value = 31251012
$ws.Range("L116").Value = 1947
$ws.Range("M116").Value = -31248718
$ws.Range("N116").Value = -6535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10001648
$ws.Range("I3").Value = 31251012
$ws.Range("J3").Value = 1947
$ws.Range("K3").Value = 31251012
$ws.Range("L3").Value = 1947
$ws.Range("M3").Value = -31250898
$ws.Range("N3").Value = -2175

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2043287.1
$ws.Range("I31").Value = 3573120
$ws.Range("J31").Value = 3510
$ws.Range("K31").Value = 3573120
$ws.Range("L31").Value = 3510
$ws.Range("M31").Value = -3572825
$ws.Range("N31").Value = -4100

$ws.Range("H34").Value = 2043287.1
$ws.Range("I34").Value = 3573120
$ws.Range("J34").Value = 3510
$ws.Range("K34").Value = 3573120
$ws.Range("L34").Value = 3510
$ws.Range("M34").Value = -3572918
$ws.Range("N34").Value = -3914

$ws.Range("H50").Value = 32500
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 32500
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 32500
$ws.Range("N50").Value = -33750

$ws.Range("H51").Value = 19274.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19274.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19274.75
$ws.Range("N51").Value = -20746.75

$ws.Range("H61").Value = 19274.75
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 19274.75
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 19274.75
$ws.Range("N61").Value = -19970.75

$ws.Range("H132").Value = 3697.8708
$ws.Range("I132").Value = 2595.4285
$ws.Range("J132").Value = 4605.7646
$ws.Range("K132").Value = 7786.2855
$ws.Range("L132").Value = 13817.2938
$ws.Range("M132").Value = -5256.2855
$ws.Range("N132").Value = -18877.2938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 653.3333
$ws.Range("I17").Value = 260
$ws.Range("J17").Value = 850
$ws.Range("K17").Value = 780
$ws.Range("L17").Value = 2550
$ws.Range("M17").Value = -611
$ws.Range("N17").Value = -2888

$ws.Range("H68").Value = 1813.5205
$ws.Range("I68").Value = 696.5185
$ws.Range("J68").Value = 2469.152
$ws.Range("K68").Value = 2089.5555
$ws.Range("L68").Value = 7407.456
$ws.Range("M68").Value = -1278.5555
$ws.Range("N68").Value = -9029.456

$ws.Range("H71").Value = 1813.5205
$ws.Range("I71").Value = 696.5185
$ws.Range("J71").Value = 2469.152
$ws.Range("K71").Value = 6268.6665
$ws.Range("L71").Value = 22222.368
$ws.Range("M71").Value = -2212.6665
$ws.Range("N71").Value = -30334.368

$ws.Range("H131").Value = 1361.7273
$ws.Range("I131").Value = 1539.5
$ws.Range("J131").Value = 1213.5834
$ws.Range("K131").Value = 4618.5
$ws.Range("L131").Value = 3640.7502
$ws.Range("M131").Value = 421.5
$ws.Range("N131").Value = -13720.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 3322.4546
$ws.Range("I132").Value = 2443.4348
$ws.Range("J132").Value = 5344.2
$ws.Range("K132").Value = 7330.3044
$ws.Range("L132").Value = 16032.6
$ws.Range("M132").Value = -4800.3044
$ws.Range("N132").Value = -21092.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2787.4285
$ws.Range("I7").Value = 1585.3334
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 1585.3334
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -1473.3334
$ws.Range("N7").Value = -10224

$ws.Range("H46").Value = 1875
$ws.Range("I46").Value = 890
$ws.Range("J46").Value = 2367.5
$ws.Range("K46").Value = 890
$ws.Range("L46").Value = 2367.5
$ws.Range("M46").Value = -702
$ws.Range("N46").Value = -2743.5

$ws.Range("H126").Value = 2787.4285
$ws.Range("I126").Value = 1585.3334
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 4756.0002
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -2286.0002
$ws.Range("N126").Value = -34940

$ws.Range("H132").Value = 5644.222
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 5971.143
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 17913.429
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -22973.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3284.6924
$ws.Range("I62").Value = 3062.25
$ws.Range("J62").Value = 3640.6
$ws.Range("K62").Value = 3062.25
$ws.Range("L62").Value = 3640.6
$ws.Range("M62").Value = -2438.25
$ws.Range("N62").Value = -4888.6

$ws.Range("H65").Value = 3284.6924
$ws.Range("I65").Value = 3062.25
$ws.Range("J65").Value = 3640.6
$ws.Range("K65").Value = 15311.25
$ws.Range("L65").Value = 18203
$ws.Range("M65").Value = -12191.25
$ws.Range("N65").Value = -24443

$ws.Range("H132").Value = 442474
$ws.Range("I132").Value = 771877.25
$ws.Range("J132").Value = 14249.8
$ws.Range("K132").Value = 2315631.75
$ws.Range("L132").Value = 42749.39999999999
$ws.Range("M132").Value = -2313101.75
$ws.Range("N132").Value = -47809.39999999999

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
